$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clean up the old placeholder content:
#    - B5 had the throw-away word "asdasdas" under the example block.
#    - A13 had the "**** COMPLETAR ****" placeholder that told the student
#      to fill in their own structure there.
# ---------------------------------------------------------------------------
$ws.Range("B5").ClearContents()
$ws.Range("A13").ClearContents()

# The row that held "asdasdas" keeps a lone centered, unbordered cell (C5)
# once the text is gone.
$ws.Range("C5").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 2. Build the real "$coleccionPalabras" structure starting at row 13,
#    mirroring the example shown above (rows 3-4) but with 10 indices
#    (columns B..K) holding the actual word list used by the Wordix game.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "`$coleccionPalabras="

function Set-IndexCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = $val
    $c.HorizontalAlignment = -4108
    $c.Font.Bold = $true
    $c.Font.Color = 10921638
}

function Set-WordCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = $val
    $c.BorderAround(1, 2, -4105, 0)
    $c.HorizontalAlignment = -4108
}

# Row 13: indices 0..9
Set-IndexCell "B13" 0
Set-IndexCell "C13" 1
Set-IndexCell "D13" 2
Set-IndexCell "E13" 3
Set-IndexCell "F13" 4
Set-IndexCell "G13" 5
Set-IndexCell "H13" 6
Set-IndexCell "I13" 7
Set-IndexCell "J13" 8
Set-IndexCell "K13" 9

# Row 14: first ten words
Set-WordCell "B14" "MUJER"
Set-WordCell "C14" "QUESO "
Set-WordCell "D14" "FUEGO"
Set-WordCell "E14" "CASAS"
Set-WordCell "F14" "RASGO"
Set-WordCell "G14" "GATOS"
Set-WordCell "H14" "GOTAS"
Set-WordCell "I14" "HUEVO"
Set-WordCell "J14" "TINTO"
Set-WordCell "K14" "NAVES"

# Row 15: indices 10..19
Set-IndexCell "B15" 10
Set-IndexCell "C15" 11
Set-IndexCell "D15" 12
Set-IndexCell "E15" 13
Set-IndexCell "F15" 14
Set-IndexCell "G15" 15
Set-IndexCell "H15" 16
Set-IndexCell "I15" 17
Set-IndexCell "J15" 18
Set-IndexCell "K15" 19

# Row 16: second ten words (F16 and K16 were filled in last, after the rest)
Set-WordCell "B16" "VERDE"
Set-WordCell "C16" "MELON"
Set-WordCell "D16" "YUYOS"
Set-WordCell "E16" "PIANO"
Set-WordCell "G16" "SILLA"
Set-WordCell "H16" "LAPIZ"
Set-WordCell "I16" "BRUMA"
Set-WordCell "J16" "RATON"
Set-WordCell "F16" "PISOS"
Set-WordCell "K16" "CABLE"

# ---------------------------------------------------------------------------
# 3. Leave the selection where Excel would after typing the last cell of the
#    new table and moving one column to the right.
# ---------------------------------------------------------------------------
$ws.Range("M16").Select()
